$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 153.69565
$ws.Range("J2").Value = 290.5
$ws.Range("L2").Value = 290.5
$ws.Range("N2").Value = -516.5
# Row 40
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4825
# Row 75
$ws.Range("H75").Value = 63078.5
$ws.Range("J75").Value = 63078.5
$ws.Range("L75").Value = 63078.5
$ws.Range("N75").Value = -64950.5
# Row 78
$ws.Range("H78").Value = 63078.5
$ws.Range("J78").Value = 63078.5
$ws.Range("L78").Value = 189235.5
$ws.Range("N78").Value = -198595.5
# Row 87
$ws.Range("H87").Value = 59999
$ws.Range("J87").Value = 59999
$ws.Range("L87").Value = 59999
$ws.Range("N87").Value = -62495
# Row 90
$ws.Range("H90").Value = 59999
$ws.Range("J90").Value = 59999
$ws.Range("L90").Value = 179997
$ws.Range("N90").Value = -192477
# Row 132
$ws.Range("H132").Value = 4153.6597
$ws.Range("I132").Value = 2895.4688
$ws.Range("K132").Value = 8686.4064
$ws.Range("M132").Value = -6156.4064

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 3152.7273
$ws.Range("I25").Value = 740.2857
$ws.Range("J25").Value = 7374.5
$ws.Range("K25").Value = 740.2857
$ws.Range("L25").Value = 7374.5
$ws.Range("M25").Value = -338.2857
$ws.Range("N25").Value = -8178.5
# Row 32
$ws.Range("H32").Value = 4378.328
$ws.Range("I32").Value = 2847.9473
$ws.Range("K32").Value = 2847.9473
$ws.Range("M32").Value = -2560.9473
# Row 61
$ws.Range("H61").Value = 3549
$ws.Range("I61").Value = 2389.8572
$ws.Range("K61").Value = 2389.8572
$ws.Range("M61").Value = -2177.8572
# Row 88
$ws.Range("H88").Value = 2438.6316
$ws.Range("J88").Value = 2404.1333
$ws.Range("L88").Value = 2404.1333
$ws.Range("N88").Value = -3216.1333
# Row 91
$ws.Range("H91").Value = 2438.6316
$ws.Range("J91").Value = 2404.1333
$ws.Range("L91").Value = 2404.1333
$ws.Range("N91").Value = -5212.1333
# Row 102
$ws.Range("H102").Value = 4498.7393
$ws.Range("I102").Value = 3269.7058
$ws.Range("K102").Value = 3269.7058
$ws.Range("M102").Value = -1647.7058
# Row 132
$ws.Range("H132").Value = 2416.1875
$ws.Range("I132").Value = 2111.0908
$ws.Range("K132").Value = 6333.2724
$ws.Range("M132").Value = -3803.2724
# Row 136
$ws.Range("H136").Value = 3549
$ws.Range("I136").Value = 2389.8572
$ws.Range("K136").Value = 7169.571599999999
$ws.Range("M136").Value = -4619.571599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2158.9656
$ws.Range("I86").Value = 2104.5715
$ws.Range("K86").Value = 2104.5715
$ws.Range("M86").Value = -981.5715
# Row 89
$ws.Range("H89").Value = 2158.9656
$ws.Range("I89").Value = 2104.5715
$ws.Range("K89").Value = 10522.8575
$ws.Range("M89").Value = -4906.8575
# Row 134
$ws.Range("H134").Value = 2586.3489
$ws.Range("I134").Value = 1898.2285
$ws.Range("K134").Value = 5694.6855
$ws.Range("M134").Value = -3159.6855
# Row 140
$ws.Range("H140").Value = 129999
$ws.Range("J140").Value = 129999
$ws.Range("L140").Value = 129999
$ws.Range("N140").Value = -140359

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2669.2173
$ws.Range("I31").Value = 2532.875
$ws.Range("J31").Value = 2980.8572
$ws.Range("K31").Value = 2532.875
$ws.Range("L31").Value = 2980.8572
$ws.Range("M31").Value = -2237.875
$ws.Range("N31").Value = -3570.8572
# Row 34
$ws.Range("H34").Value = 2669.2173
$ws.Range("I34").Value = 2532.875
$ws.Range("J34").Value = 2980.8572
$ws.Range("K34").Value = 2532.875
$ws.Range("L34").Value = 2980.8572
$ws.Range("M34").Value = -2330.875
$ws.Range("N34").Value = -3384.8572

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 1466.0526
$ws.Range("I11").Value = 410.75
$ws.Range("J11").Value = 3275.1428
$ws.Range("K11").Value = 1232.25
$ws.Range("L11").Value = 9825.428400000001
$ws.Range("M11").Value = -1092.25
$ws.Range("N11").Value = -10105.4284
# Row 29
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 46
$ws.Range("H46").Value = 259444720
$ws.Range("I46").Value = 190
$ws.Range("J46").Value = 273099680
$ws.Range("K46").Value = 570
$ws.Range("L46").Value = 819299040
$ws.Range("M46").Value = -479
$ws.Range("N46").Value = -819299222
# Row 50
$ws.Range("H50").Value = 3381.5557
$ws.Range("I50").Value = 2149
$ws.Range("K50").Value = 6447
$ws.Range("M50").Value = -5966
# Row 53
$ws.Range("H53").Value = 3381.5557
$ws.Range("I53").Value = 2149
$ws.Range("K53").Value = 6447
$ws.Range("M53").Value = -5966
# Row 98
$ws.Range("H98").Value = 1915.5
$ws.Range("I98").Value = 1999
$ws.Range("J98").Value = 1898.8
$ws.Range("K98").Value = 5997
$ws.Range("L98").Value = 5696.4
$ws.Range("M98").Value = -4499
$ws.Range("N98").Value = -8692.4
# Row 139
$ws.Range("H139").Value = 3476.0667
$ws.Range("J139").Value = 9248.75
$ws.Range("L139").Value = 27746.25
$ws.Range("N139").Value = -38026.25

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 650.2105
$ws.Range("I97").Value = 604.75
$ws.Range("K97").Value = 604.75
$ws.Range("M97").Value = -108.75
# Row 113
$ws.Range("H113").Value = 11172.25
$ws.Range("I113").Value = 3540
$ws.Range("J113").Value = 21857.4
$ws.Range("K113").Value = 3540
$ws.Range("L113").Value = 21857.4
$ws.Range("M113").Value = -1370
$ws.Range("N113").Value = -26197.4
# Row 126
$ws.Range("H126").Value = 5258.4736
$ws.Range("I126").Value = 4811.923
$ws.Range("J126").Value = 6226
$ws.Range("K126").Value = 14435.769
$ws.Range("L126").Value = 18678
$ws.Range("M126").Value = -11965.769
$ws.Range("N126").Value = -23618
# Row 132
$ws.Range("H132").Value = 2269.16
$ws.Range("I132").Value = 1666.55
$ws.Range("K132").Value = 4999.65
$ws.Range("M132").Value = -2469.65

$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 50007
$ws.Range("J13").Value = 50007
$ws.Range("L13").Value = 50007
$ws.Range("N13").Value = -50287
# Row 22
$ws.Range("H22").Value = 815.8570999999999
$ws.Range("I22").Value = 815.8570999999999
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 815.8570999999999
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -520.8570999999999
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 815.8570999999999
$ws.Range("I27").Value = 815.8570999999999
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 815.8570999999999
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -708.8570999999999
$ws.Range("N27").ClearContents()
# Row 46
$ws.Range("I46").Value = 1100
$ws.Range("J46").Value = 2563.5833
$ws.Range("K46").Value = 1100
$ws.Range("L46").Value = 2563.5833
$ws.Range("M46").Value = -912
$ws.Range("N46").Value = -2939.5833
# Row 55
$ws.Range("H55").Value = 1336.8235
$ws.Range("I55").Value = 525.2857
$ws.Range("J55").Value = 1904.9
$ws.Range("K55").Value = 525.2857
$ws.Range("L55").Value = 1904.9
$ws.Range("M55").Value = -352.2857
$ws.Range("N55").Value = -2250.9
# Row 100
$ws.Range("H100").Value = 3066.6667

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4084.3262
$ws.Range("I126").Value = 3761.2896
$ws.Range("K126").Value = 11283.8688
$ws.Range("M126").Value = -8813.8688
# Row 132
$ws.Range("H132").Value = 1788.6666
$ws.Range("I132").Value = 1577
$ws.Range("J132").Value = 1923.3636
$ws.Range("K132").Value = 4731
$ws.Range("L132").Value = 5770.0908
$ws.Range("M132").Value = -2201
$ws.Range("N132").Value = -10830.0908

